$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1210.5
$ws.Range("I20").Value = 1210.5
$ws.Range("K20").Value = 1210.5
$ws.Range("M20").Value = -980.5

$ws.Range("H35").Value = 1210.5
$ws.Range("I35").Value = 1210.5
$ws.Range("K35").Value = 1210.5
$ws.Range("M35").Value = -831.5

$ws.Range("H94").Value = 6505.263
$ws.Range("I94").Value = 3966.6667
$ws.Range("J94").Value = 7676.923
$ws.Range("K94").Value = 3966.6667
$ws.Range("L94").Value = 7676.923
$ws.Range("M94").Value = -3515.6667
$ws.Range("N94").Value = -8578.922999999999

$ws.Range("H96").Value = 453.75
$ws.Range("I96").Value = 362.66666
$ws.Range("J96").Value = 544.8333
$ws.Range("K96").Value = 1087.99998
$ws.Range("L96").Value = 1634.4999
$ws.Range("M96").Value = 285.0000199999999
$ws.Range("N96").Value = -4380.4999

$ws.Range("H99").Value = 648
$ws.Range("I99").Value = 222
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 666
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = 832
$ws.Range("N99").Value = -7496

$ws.Range("H100").Value = 83334970
$ws.Range("I100").Value = 1961
$ws.Range("J100").Value = 500000000
$ws.Range("K100").Value = 1961
$ws.Range("L100").Value = 500000000
$ws.Range("M100").Value = -1420
$ws.Range("N100").Value = -500001082

$ws.Range("H112").Value = 3726.5715
$ws.Range("I112").Value = 1033.3334
$ws.Range("J112").Value = 3979.0625
$ws.Range("K112").Value = 3100.0002
$ws.Range("L112").Value = 11937.1875
$ws.Range("M112").Value = -1992.0002
$ws.Range("N112").Value = -14153.1875

$ws.Range("H113").Value = 3267.0435
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 4979.2
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 4979.2
$ws.Range("M113").Value = 1304
$ws.Range("N113").Value = -11487.2

$ws.Range("H136").Value = 85853.336
$ws.Range("J136").Value = 85853.336
$ws.Range("L136").Value = 85853.336
$ws.Range("N136").Value = -96053.336

$ws.Range("H141").Value = 3764.1428
$ws.Range("I141").Value = 3424.8333
$ws.Range("K141").Value = 10274.4999
$ws.Range("M141").Value = -5094.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 16814.572
$ws.Range("J55").Value = 16814.572
$ws.Range("L55").Value = 16814.572
$ws.Range("N55").Value = -17444.572

$ws.Range("H80").Value = 24970.8
$ws.Range("J80").Value = 24970.8
$ws.Range("L80").Value = 24970.8
$ws.Range("N80").Value = -26966.8

$ws.Range("H83").Value = 24970.8
$ws.Range("J83").Value = 24970.8
$ws.Range("L83").Value = 74912.39999999999
$ws.Range("N83").Value = -84896.39999999999

$ws.Range("H102").Value = 142857980
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 500000450
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 500000450
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -500003694

$ws.Range("H113").Value = 43776
$ws.Range("J113").Value = 43776
$ws.Range("L113").Value = 43776
$ws.Range("N113").Value = -52454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27832.15
$ws.Range("J82").Value = 33617.734
$ws.Range("L82").Value = 33617.734
$ws.Range("N82").Value = -34383.734

$ws.Range("H85").Value = 27832.15
$ws.Range("J85").Value = 33617.734
$ws.Range("L85").Value = 33617.734
$ws.Range("N85").Value = -36269.734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H41").Value = 17062.5
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 20375
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 20375
$ws.Range("M41").Value = -72
$ws.Range("N41").Value = -21231

$ws.Range("H50").Value = 9019.714
$ws.Range("J50").Value = 9019.714
$ws.Range("L50").Value = 9019.714
$ws.Range("N50").Value = -10269.714

$ws.Range("H51").Value = 9232
$ws.Range("J51").Value = 9232
$ws.Range("L51").Value = 9232
$ws.Range("N51").Value = -10704

$ws.Range("H58").Value = 1746.7858
$ws.Range("I58").Value = 1746.7858
$ws.Range("K58").Value = 1746.7858
$ws.Range("M58").Value = -1543.7858

$ws.Range("H60").Value = 30509.334
$ws.Range("J60").Value = 30509.334
$ws.Range("L60").Value = 30509.334
$ws.Range("N60").Value = -31531.334

$ws.Range("H61").Value = 9232
$ws.Range("J61").Value = 9232
$ws.Range("L61").Value = 9232
$ws.Range("N61").Value = -9928

$ws.Range("H68").Value = 17466.334
$ws.Range("J68").Value = 17466.334
$ws.Range("L68").Value = 17466.334
$ws.Range("N68").Value = -18964.334

$ws.Range("H71").Value = 17466.334
$ws.Range("J71").Value = 17466.334
$ws.Range("L71").Value = 52399.00199999999
$ws.Range("N71").Value = -59887.00199999999

$ws.Range("H109").Value = 10549.375
$ws.Range("J109").Value = 10549.375
$ws.Range("L109").Value = 10549.375
$ws.Range("N109").Value = -12629.375

$ws.Range("H122").Value = 1765.4286
$ws.Range("I122").Value = 1048.1428
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 3144.4284
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -694.4284000000002
$ws.Range("N122").Value = -14500

$ws.Range("H134").Value = 8720.143
$ws.Range("I134").Value = 9269.632
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 27808.896
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -25273.896
$ws.Range("N134").Value = -15570

$ws.Range("H136").Value = 1746.7858
$ws.Range("I136").Value = 1746.7858
$ws.Range("K136").Value = 5240.357400000001
$ws.Range("M136").Value = -2690.357400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 76.25
$ws.Range("I2").Value = 76.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 76.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 36.75
$ws.Range("N2").ClearContents()

$ws.Range("H57").Value = 18030.428
$ws.Range("J57").Value = 20202.334
$ws.Range("L57").Value = 20202.334
$ws.Range("N57").Value = -21842.334

$ws.Range("H123").Value = 33605
$ws.Range("J123").Value = 33605
$ws.Range("L123").Value = 33605
$ws.Range("N123").Value = -38505

$ws.Range("H136").Value = 19841.2
$ws.Range("J136").Value = 19841.2
$ws.Range("L136").Value = 59523.60000000001
$ws.Range("N136").Value = -64623.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1610
$ws.Range("I100").Value = 1610
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1610
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1069
$ws.Range("N100").ClearContents()

$ws.Range("H136").Value = 1870.1708
$ws.Range("I136").Value = 1416.9706
$ws.Range("J136").Value = 4071.4285
$ws.Range("K136").Value = 4250.9118
$ws.Range("L136").Value = 12214.2855
$ws.Range("M136").Value = -1700.9118
$ws.Range("N136").Value = -17314.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 100000000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H109").Value = 16555.666
$ws.Range("J109").Value = 16555.666
$ws.Range("L109").Value = 16555.666
$ws.Range("N109").Value = -19329.666

$ws.Range("H136").Value = 2486.6086
$ws.Range("I136").Value = 2090.6365
$ws.Range("J136").Value = 2849.5833
$ws.Range("K136").Value = 6271.9095
$ws.Range("L136").Value = 8548.749899999999
$ws.Range("M136").Value = -3721.9095
$ws.Range("N136").Value = -13648.7499
